# "final changes for diss"
# - Zoom the active sheet view to 112%
# - Change the selection to A2:A38 (active cell A2)
# - Widen column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# Select the worksheet's data range A2:A38, with A2 as the active cell
$ws.Range("A2:A38").Select()

# Set the zoom level for the view to 112%
$win.Zoom = 112

# Widen column B (from ~13.86 to ~18.14 characters)
$ws.Columns.Item(2).ColumnWidth = 17.333333333333332
